# Insert a new data row at row 316 (pushing the existing rows 316-383
# down to 317-384) and populate it with the new observation.
# This mirrors the diff: dimension grows from A1:R383 to A1:R384, and
# every row from the old 316 onward is shifted down by one; the freshly
# inserted row 316 carries brand-new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(316).Insert()

$ws.Cells.Item(316, 1).Value = 10
$ws.Cells.Item(316, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(316, 3).Value = "La Araucanía"
$ws.Cells.Item(316, 4).Value = 44754
$ws.Cells.Item(316, 5).Value = 9
$ws.Cells.Item(316, 6).Value = 100112037
$ws.Cells.Item(316, 7).Value = "Cebollín"
$ws.Cells.Item(316, 8).Value = "Sin especificar"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 40
$ws.Cells.Item(316, 11).Value = 7000
$ws.Cells.Item(316, 12).Value = 7000
$ws.Cells.Item(316, 13).Value = 7000
$ws.Cells.Item(316, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(316, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(316, 16).Value = 583
$ws.Cells.Item(316, 17).Value = 12
$ws.Cells.Item(316, 18).Value = "Hortaliza"

# Match the date-format style used by the rest of column D (style index 2
# in the original workbook corresponds to the yyyy-mm-dd style applied to
# the "Fecha" column).
$ws.Cells.Item(316, 4).NumberFormat = $ws.Cells.Item(317, 4).NumberFormat
